$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-numeric-looking text values (plain assignment is safe)
$ws.Range("D2").Value = '42.839.74'
$ws.Range("E2").Value = '  +4.14%  '
$ws.Range("D3").Value = '2.254.39'
$ws.Range("E3").Value = '  +3.19%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("E6").Value = '  +1.76%  '
$ws.Range("E7").Value = '  +4.47%  '
$ws.Range("E8").Value = '  +17.93%  '
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("E10").Value = '  +10.54%  '
$ws.Range("E11").Value = '  +4.10%  '
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("E13").Value = '  +8.51%  '
$ws.Range("E14").Value = '  +1.11%  '
$ws.Range("D15").Value = '2.592.28'
$ws.Range("E15").Value = '  +2.93%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("E16").Value = '  +4.55%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("E17").Value = '  +2.08%  '
$ws.Range("D18").Value = '2.266.17'
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("D19").Value = '42.736.42'
$ws.Range("E19").Value = '  +4.06%  '
$ws.Range("E20").Value = '  +3.47%  '
$ws.Range("E21").Value = '  +2.30%  '
$ws.Range("E22").Value = '  +1.66%  '
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("E24").Value = '  -1.77%  '
$ws.Range("E25").Value = '  +6.18%  '
$ws.Range("E26").Value = '  +2.06%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("E30").Value = '  +3.02%  '
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("E32").Value = '  +3.08%  '
$ws.Range("E33").Value = '  +12.98%  '
$ws.Range("E34").Value = '  +7.04%  '
$ws.Range("E35").Value = '  +3.96%  '
$ws.Range("E36").Value = '  +20.22%  '
$ws.Range("E37").Value = '  +2.98%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E38").Value = '  +1.45%  '
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("E39").Value = '  +2.40%  '
$ws.Range("E40").Value = '  +4.83%  '
$ws.Range("E41").Value = '  +5.17%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("E42").Value = '  +5.56%  '
$ws.Range("B43").Value = 'THORChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("E43").Value = '  +4.98%  '
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("E44").Value = '  +1.15%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("E46").Value = '  +3.10%  '
$ws.Range("E47").Value = '  +4.11%  '
$ws.Range("E48").Value = '  +4.18%  '
$ws.Range("E50").Value = '  -3.49%  '
$ws.Range("E51").Value = '  +3.09%  '

# Numeric-looking values that must remain text -> force text format, assign, then reset style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '251.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.634'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '70.41'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.677'
$ws.Range("D8").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0974'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '59.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.90'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.889'
$ws.Range("D17").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.78'
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.46'
$ws.Range("D28").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.04'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.34'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.128'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0787'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '30.87'
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.34'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.71'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0322'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.32'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.86'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.92'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.202'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.20'
$ws.Range("D51").Style = "Normal"